# Scheduled runner update: refresh market-price-derived Leve profit figures
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ columns)
# across the per-class Leve tables.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 21828.5
$ws.Range("J3").Value = 21828.5
$ws.Range("L3").Value = 21828.5
$ws.Range("N3").Value = -22056.5
$ws.Range("H4").Value = 997.5
$ws.Range("I4").Value = 997.5
$ws.Range("K4").Value = 997.5
$ws.Range("M4").Value = -883.5
$ws.Range("H11").Value = 63
$ws.Range("I11").Value = 63
$ws.Range("K11").Value = 63
$ws.Range("M11").Value = 77
$ws.Range("H102").Value = 21828.5
$ws.Range("J102").Value = 21828.5
$ws.Range("L102").Value = 21828.5
$ws.Range("N102").Value = -28318.5
$ws.Range("H138").Value = 2781.7778
$ws.Range("I138").Value = 1811.125
$ws.Range("K138").Value = 5433.375
$ws.Range("M138").Value = -293.375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 225459.89
$ws.Range("I32").Value = 2276.8462
$ws.Range("K32").Value = 2276.8462
$ws.Range("M32").Value = -1989.8462
$ws.Range("H45").Value = 2738.9
$ws.Range("I45").Value = 2148.6667
$ws.Range("J45").Value = 3624.25
$ws.Range("K45").Value = 2148.6667
$ws.Range("L45").Value = 3624.25
$ws.Range("M45").Value = -1771.6667
$ws.Range("N45").Value = -4378.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 433.9091
$ws.Range("I11").Value = 86
$ws.Range("J11").Value = 564.375
$ws.Range("K11").Value = 86
$ws.Range("L11").Value = 564.375
$ws.Range("M11").Value = 54
$ws.Range("N11").Value = -844.375
$ws.Range("H94").Value = 1126.1428
$ws.Range("I94").Value = 576.6
$ws.Range("K94").Value = 576.6
$ws.Range("M94").Value = -125.6

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1226.2
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -326
$ws.Range("H86").Value = 8217.200000000001
$ws.Range("I86").Value = 7832.6665
$ws.Range("K86").Value = 7832.6665
$ws.Range("M86").Value = -6709.6665
$ws.Range("H89").Value = 8217.200000000001
$ws.Range("I89").Value = 7832.6665
$ws.Range("K89").Value = 39163.3325
$ws.Range("M89").Value = -33547.3325
$ws.Range("H132").Value = 4564.125
$ws.Range("J132").Value = 6099.1665
$ws.Range("L132").Value = 18297.4995
$ws.Range("N132").Value = -23357.4995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2000
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 6000
$ws.Range("N5").Value = -6224
$ws.Range("H92").Value = 537.75
$ws.Range("J92").Value = 477.5
$ws.Range("L92").Value = 1432.5
$ws.Range("N92").Value = -3928.5
$ws.Range("H135").Value = 2000
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("H25").Value = 576.1875
$ws.Range("I25").Value = 486.85715
$ws.Range("J25").Value = 645.6667
$ws.Range("K25").Value = 486.85715
$ws.Range("L25").Value = 645.6667
$ws.Range("M25").Value = 42.14285000000001
$ws.Range("N25").Value = -1703.6667
$ws.Range("H44").Value = 13816.667
$ws.Range("J44").Value = 9450
$ws.Range("L44").Value = 9450
$ws.Range("N44").Value = -10642
$ws.Range("H97").Value = 684.1667
$ws.Range("I97").Value = 532.1875
$ws.Range("K97").Value = 532.1875
$ws.Range("M97").Value = -36.1875
$ws.Range("H132").Value = 4025.25
$ws.Range("I132").Value = 3608.7144
$ws.Range("J132").Value = 4608.4
$ws.Range("K132").Value = 10826.1432
$ws.Range("L132").Value = 13825.2
$ws.Range("M132").Value = -8296.143199999999
$ws.Range("N132").Value = -18885.2
$ws.Range("M5").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1327.7142
$ws.Range("I22").Value = 1149.75
$ws.Range("J22").Value = 1565
$ws.Range("K22").Value = 1149.75
$ws.Range("L22").Value = 1565
$ws.Range("M22").Value = -854.75
$ws.Range("N22").Value = -2155
$ws.Range("H27").Value = 1327.7142
$ws.Range("I27").Value = 1149.75
$ws.Range("J27").Value = 1565
$ws.Range("K27").Value = 1149.75
$ws.Range("L27").Value = 1565
$ws.Range("M27").Value = -1042.75
$ws.Range("N27").Value = -1779
$ws.Range("H46").Value = 2333.1333
$ws.Range("I46").Value = 1025.25
$ws.Range("J46").Value = 3827.8572
$ws.Range("K46").Value = 1025.25
$ws.Range("L46").Value = 3827.8572
$ws.Range("M46").Value = -837.25
$ws.Range("N46").Value = -4203.8572
$ws.Range("H55").Value = 1145.2727
$ws.Range("I55").Value = 1024.875
$ws.Range("J55").Value = 1466.3334
$ws.Range("K55").Value = 1024.875
$ws.Range("L55").Value = 1466.3334
$ws.Range("M55").Value = -851.875
$ws.Range("N55").Value = -1812.3334
$ws.Range("H61").Value = 1633.5
$ws.Range("I61").Value = 1705.8334
$ws.Range("J61").Value = 1199.5
$ws.Range("K61").Value = 1705.8334
$ws.Range("L61").Value = 1199.5
$ws.Range("M61").Value = -1503.8334
$ws.Range("N61").Value = -1603.5
$ws.Range("H113").Value = 1633.5
$ws.Range("I113").Value = 1705.8334
$ws.Range("J113").Value = 1199.5
$ws.Range("K113").Value = 1705.8334
$ws.Range("L113").Value = 1199.5
$ws.Range("M113").Value = 464.1666
$ws.Range("N113").Value = -5539.5
$ws.Range("H134").Value = 33814.5
$ws.Range("J134").Value = 33814.5
$ws.Range("L134").Value = 33814.5
$ws.Range("N134").Value = -43954.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 29498
$ws.Range("J33").Value = 27997
$ws.Range("L33").Value = 27997
$ws.Range("N33").Value = -28497
$ws.Range("H36").Value = 29498
$ws.Range("J36").Value = 27997
$ws.Range("L36").Value = 27997
$ws.Range("N36").Value = -28497
$ws.Range("H45").Value = 16064
$ws.Range("I45").Value = 9899
$ws.Range("J45").Value = 17605.25
$ws.Range("K45").Value = 9899
$ws.Range("L45").Value = 17605.25
$ws.Range("M45").Value = -9408
$ws.Range("N45").Value = -18587.25
$ws.Range("H103").Value = 13468
$ws.Range("J103").Value = 13468
$ws.Range("L103").Value = 13468
$ws.Range("N103").Value = -15812
$ws.Range("H104").Value = 17586.625
$ws.Range("J104").Value = 17586.625
$ws.Range("L104").Value = 17586.625
$ws.Range("N104").Value = -24574.625
